$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 7 new rows before the old "subtotal" row (old row 9), shifting
# everything below down by 7 rows (old row 9 -> new row 16, etc.)
$ws.Rows("9:15").Insert()

# --- Row 9: 3.3v regulator ---
$ws.Range("C9").Value = 0.04
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = "3.3v regulator"
$ws.Range("F9").Value = "http://www.aliexpress.com/item/10pcs-AMS1117-3-3-AMS1117-LM1117-1117-3-3V-1A-Voltage-Regulator/32409097011.html"

# --- Row 10: micro usb connector ---
$ws.Range("C10").Value = 0.08
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = "micro usb connector"
$ws.Range("F10").Value = "http://www.aliexpress.com/item/Short-term-Sacrifices-5pcs-USB-Micro-Type-B-5pin-Female-Jack-Connector-SMT-Surface-Mount/32349977492.html"

# --- Row 11: micro usb cable ---
$ws.Range("C11").Value = 0.72
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = "micro usb cable"
$ws.Range("F11").Value = "http://www.aliexpress.com/item/1-Meter-3-Feet-Ruggedized-Fabric-Braided-USB-Male-to-Micro-USB-Male-Data-Sync-Charging/32340270515.html"

# --- Row 12: tft screen ---
$ws.Range("C12").Value = 3.6
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = "tft screen"
$ws.Range("F12").Value = "http://www.aliexpress.com/item/Free-shipping-New-1-44-inch-LCD-color-screen-for-arduino-1-44-TFT-SPI-serial/32233128819.html"

# --- Row 13: rgb led ---
$ws.Range("C13").Value = 0.75
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = "rgb led"
$ws.Range("F13").Value = "http://www.aliexpress.com/item/1Pcs-3-Colour-RGB-SMD-LED-Board-Module-5050-Full-Color-LED-3-3-5V-for/2055617889.html"

# --- Row 14: circuit board (no link) ---
$ws.Range("C14").Value = 5
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = "circuit board"

# --- Row 15: case (no link) ---
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = "case"

# Fill the B column (cost = price * qty) formula down through the new rows
$ws.Range("B9").Formula = "=C9*D9"
$ws.Range("B10").Formula = "=C10*D10"
$ws.Range("B11").Formula = "=C11*D11"
$ws.Range("B12").Formula = "=C12*D12"
$ws.Range("B13").Formula = "=C13*D13"
$ws.Range("B14").Formula = "=C14*D14"
$ws.Range("B15").Formula = "=C15*D15"

# Extend the subtotal SUM to cover the new rows
$ws.Range("B16").Formula = "=SUM(B2:B15)"

# Update the selected cell to match the author's final selection
[void]$ws.Range("F15").Select()
